$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "datepicker selecter"
$ws.Range("B21").Value = "find by xpath"
$ws.Range("C21").Value = "//button[@id='ext-gen23']"

$ws.Range("B22").Value = "find by xpath"
$ws.Range("C22").Value = "(//td[contains(@class,'x-date-active')]//span)[@@]"
$ws.Range("D22").Value = "getcurrentdate"
